$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 192, shifting existing rows 192-229 down to 196-233
$ws.Range("A192:A195").EntireRow.Insert()

# Row 192
$ws.Cells.Item(192, 1).Value = 3
$ws.Cells.Item(192, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(192, 3).Value = 'Coquimbo'
$ws.Cells.Item(192, 4).Value = 44522
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 'Fruta'
$ws.Cells.Item(192, 7).Value = 100103
$ws.Cells.Item(192, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(192, 9).Value = 100103001
$ws.Cells.Item(192, 10).Value = 'Cereza'
$ws.Cells.Item(192, 11).Value = 'Royal Dawn'
$ws.Cells.Item(192, 12).Value = 'Primera'
$ws.Cells.Item(192, 13).Value = 165
$ws.Cells.Item(192, 14).Value = 22000
$ws.Cells.Item(192, 15).Value = 23000
$ws.Cells.Item(192, 16).Value = 22485
$ws.Cells.Item(192, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(192, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(192, 19).Value = 2248
$ws.Cells.Item(192, 20).Value = 10

# Row 193
$ws.Cells.Item(193, 1).Value = 3
$ws.Cells.Item(193, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(193, 3).Value = 'Coquimbo'
$ws.Cells.Item(193, 4).Value = 44522
$ws.Cells.Item(193, 5).Value = 5
$ws.Cells.Item(193, 6).Value = 'Fruta'
$ws.Cells.Item(193, 7).Value = 100103
$ws.Cells.Item(193, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(193, 9).Value = 100103001
$ws.Cells.Item(193, 10).Value = 'Cereza'
$ws.Cells.Item(193, 11).Value = 'Royal Dawn'
$ws.Cells.Item(193, 12).Value = 'Segunda'
$ws.Cells.Item(193, 13).Value = 80
$ws.Cells.Item(193, 14).Value = 18000
$ws.Cells.Item(193, 15).Value = 18000
$ws.Cells.Item(193, 16).Value = 18000
$ws.Cells.Item(193, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(193, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(193, 19).Value = 1800
$ws.Cells.Item(193, 20).Value = 10

# Row 194
$ws.Cells.Item(194, 1).Value = 3
$ws.Cells.Item(194, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(194, 3).Value = 'Coquimbo'
$ws.Cells.Item(194, 4).Value = 44522
$ws.Cells.Item(194, 5).Value = 5
$ws.Cells.Item(194, 6).Value = 'Fruta'
$ws.Cells.Item(194, 7).Value = 100103
$ws.Cells.Item(194, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(194, 9).Value = 100103001
$ws.Cells.Item(194, 10).Value = 'Cereza'
$ws.Cells.Item(194, 11).Value = 'Santina'
$ws.Cells.Item(194, 12).Value = 'Primera'
$ws.Cells.Item(194, 13).Value = 60
$ws.Cells.Item(194, 14).Value = 23000
$ws.Cells.Item(194, 15).Value = 23000
$ws.Cells.Item(194, 16).Value = 23000
$ws.Cells.Item(194, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(194, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(194, 19).Value = 2300
$ws.Cells.Item(194, 20).Value = 10

# Row 195
$ws.Cells.Item(195, 1).Value = 3
$ws.Cells.Item(195, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(195, 3).Value = 'Coquimbo'
$ws.Cells.Item(195, 4).Value = 44522
$ws.Cells.Item(195, 5).Value = 5
$ws.Cells.Item(195, 6).Value = 'Fruta'
$ws.Cells.Item(195, 7).Value = 100103
$ws.Cells.Item(195, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(195, 9).Value = 100103001
$ws.Cells.Item(195, 10).Value = 'Cereza'
$ws.Cells.Item(195, 11).Value = 'Santina'
$ws.Cells.Item(195, 12).Value = 'Segunda'
$ws.Cells.Item(195, 13).Value = 55
$ws.Cells.Item(195, 14).Value = 19000
$ws.Cells.Item(195, 15).Value = 19000
$ws.Cells.Item(195, 16).Value = 19000
$ws.Cells.Item(195, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(195, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(195, 19).Value = 1900
$ws.Cells.Item(195, 20).Value = 10
